$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Sheet "Resumen por Regimen" (Worksheets index 1): swap F31/G31 and
#    F32/G32 values, then replace the view (drop tabSelected/topLeftCell,
#    add frozen panes at C2 i.e. 1 row / 2 cols frozen, with the bottom-right
#    pane scrolled so the selection ends on F35:F36).
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()

$f31 = $ws1.Range("F31").Value()
$g31 = $ws1.Range("G31").Value()
$ws1.Range("F31").Value = $g31
$ws1.Range("G31").Value = $f31

$f32 = $ws1.Range("F32").Value()
$g32 = $ws1.Range("G32").Value()
$ws1.Range("F32").Value = $g32
$ws1.Range("G32").Value = $f32

$ws1.Range("C2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws1.Range("F35:F36").Select()

# ---------------------------------------------------------------------------
# 2) Sheet "Resumen por Categoria" (Worksheets index 2): swap F/G for rows
#    87-95, then replace the view the same way (frozen panes, bottom-right
#    pane selection F98:G110).
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()

foreach ($r in 87..95) {
    $fCell = "F" + $r
    $gCell = "G" + $r
    $fVal = $ws2.Range($fCell).Value()
    $gVal = $ws2.Range($gCell).Value()
    $ws2.Range($fCell).Value = $gVal
    $ws2.Range($gCell).Value = $fVal
}

$ws2.Range("C2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws2.Range("F98:G110").Select()

# ---------------------------------------------------------------------------
# 3) Sheet "Theil por Categoria" (Worksheets index 5) becomes the active /
#    selected tab.
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Activate()
$ws5.Range("D1").Select()
